# ---------------------------------------------------------------------------
# Applies the edit described by the commit "blood pressure map, and
# ELIXIR-LU tryout":
#   1. Refreshes the cached "datetimeFigureOut" footer field text (found on
#      the slide master and every slide layout) from 18/12/2024 to 21/10/2025.
#   2. Adds a new "Blood pressure" rounded-rectangle legend entry shape to
#      slide 1, matching the existing star-rating legend shapes already on
#      the slide.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Step 1: update the cached date text on the slide master + all layouts.
# ---------------------------------------------------------------------------
$oldDate = "18/12/2024"
$newDate = "21/10/2025"

$master = $p.SlideMaster

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Update-DatePlaceholder $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# Step 2: add the new "Blood pressure" legend shape to slide 1.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# Locate an existing legend entry shape to clone so the new shape inherits
# the same preset geometry / fill / line / style references used by every
# other rounded-rectangle legend entry on this slide.
$template = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Rounded Rectangle 48") {
        $template = $shp
    }
}

$newShapeRange = $template.Duplicate()
$newShape = $newShapeRange.Item(1)

$newShape.Name = "Rounded Rectangle 4"

$newShape.Left = 9664435 / 12700
$newShape.Top = 4652124 / 12700
$newShape.Width = 1740988 / 12700
$newShape.Height = 581006 / 12700

$newShape.TextFrame.VerticalAnchor = 1
$newShape.TextFrame.TextRange.Text = "Blood pressure`r★☆☆☆☆"

Write-Host "Added shape id=$($newShape.Id) name=$($newShape.Name)"
